# [Fonds de solidarite] Add 2020-12-09 data
# Update nombre_aides (col C) and montant_total (col D) for the rows whose
# underlying source data changed with the new 2020-12-09 snapshot.
# The sheet stores these columns as text (inlineStr) values, so we force
# the Text number format before writing so Excel doesn't re-type the
# value as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => nombre_aides (C), montant_total (D)
$updates = @{
    3   = @("1385", "6376555.24")
    4   = @("598",  "4952932.87")
    25  = @("73",   "935024.25")
    56  = @("1010", "5588551.94")
    57  = @("486",  "4333391.40")
    58  = @("177",  "1700340.18")
    59  = @("67",   "995420.00")
    60  = @("14",   "301032.25")
    63  = @("5688", "23860987.63")
    64  = @("3108", "18466223.98")
    105 = @("495",  "2305228.10")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]

    $cCell = $ws.Cells.Item($row, 3)
    $dCell = $ws.Cells.Item($row, 4)

    $cCell.NumberFormat = "@"
    $dCell.NumberFormat = "@"

    $cCell.Value = $values[0]
    $dCell.Value = $values[1]
}
